# Adds a new worksheet "HPT3An" (third-degree/3-unknown linear-system test
# cases) after the existing "PTB2" sheet, mirroring the layout of the first
# two sheets, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet after the last existing sheet ("PTB2") and
#    rename it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "HPT3An"

# ---------------------------------------------------------------------
# 2. Header block (rows 1-2), copied verbatim from the other sheets.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Hàm"
$ws.Range("C1").Value = "Giải phương trình Bậc 2"
$ws.Range("A2").Value = "Các trường hợp test : "

$ws.Range("A1:B1").Merge()
$ws.Range("C1:G1").Merge()
$ws.Range("A2:G2").Merge()

# ---------------------------------------------------------------------
# 3. Column header rows (3 & 4).
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "TestCase"
$ws.Range("B3").Value = "Mô tả"
$ws.Range("C3").Value = "Dữ liệu nhập"
$ws.Range("D3").Value = "Kết quả mong đợi"
$ws.Range("E3").Value = "Kểt quả chạy"
$ws.Range("F3").Value = "Failed/ Passed"
$ws.Range("G3").Value = "Report `n(Nếu failed) "

$ws.Range("A4").Value = "<TestCase ID>"
$ws.Range("B4").Value = "<Mô tả về trường hợp test >"
$ws.Range("C4").Value = "<Mô tả dữ liệu nhập>"
$ws.Range("D4").Value = "<Mô tả kết quả mong đợi>"
$ws.Range("E4").Value = " <Kểt quả khi chạy chương trình>"
$ws.Range("F4").Value = "<Failed hay Passed ?>"
$ws.Range("G4").Value = "<Đưa ra báo cáo cho trường hợp bị Failed>"

# ---------------------------------------------------------------------
# 4. Test-case rows (5-7): 3x3 linear system test cases.
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "TC_HPT3A_01"
$ws.Range("B5").Value = "Kiểm tra hệ phương trình vô số nghiệm"
$ws.Range("C5").Value = "Nhập vào a1=0, b1=0, c1=0, d1=0, a2=0, b2=0, c2=0, d2=0, a3=0, b3=0, c3=0, d3=0"
$ws.Range("D5").Value = 'Thông báo "Hệ phương trình có vô số nghiệm"'
$ws.Range("E5").Value = "Hệ phương trình có vô số nghiệm"
$ws.Range("F5").Value = "Passed"

$ws.Range("A6").Value = "TC_HPT3A_02"
$ws.Range("B6").Value = "Kiểm tra hệ phương trình vô nghiệm"
$ws.Range("C6").Value = "Nhập vào a1=0, b1=0, c1=0, d1=1, a2=0, b2=0, c2=0, d2=0, a3=0, b3=0, c3=0, d3=0"
$ws.Range("D6").Value = 'Thông báo "Hệ phương trình vô nghiệm"'
$ws.Range("E6").Value = "Hệ phương trình có vô số nghiệm"
$ws.Range("F6").Value = "Failed"
$ws.Range("G6").Value = "Message:   Expected string length 25 but was 28. Strings differ at index 19.`n  Expected: `"Hệ phương trình vô nghiệm`"`n  But was:  `"Hệ phương trình vô số nghiệm`"`n  ------------------------------^"

$ws.Range("A7").Value = "TC_HPT3A_03"
$ws.Range("B7").Value = "Kiểm tra hệ phương trình có một nghiệm"
$ws.Range("C7").Value = "Nhập vào a1=1, b1=1, c1=1, d1=6, a2=1, b2=-1, c2=1, d2=2, a3=1, b3=1, c3=-1, d3=0"
$ws.Range("D7").Value = "Thông báo `"Hệ phương trình có 1 nghiệm`"`nKết quả X=1, Y=2, Z=3"
$ws.Range("E7").Value = "Hệ phương trình vô nghiệm"
$ws.Range("F7").Value = "Failed"
$ws.Range("G7").Value = "Message:   Expected string length 27 but was 25. Strings differ at index 16.`n  Expected: `"Hệ phương trình có 1 nghiệm`"`n  But was:  `"Hệ phương trình vô nghiệm`"`n  ---------------------------^"

# ---------------------------------------------------------------------
# 5. Formatting: fonts / borders / alignment, matching the TestCase-row
#    styling of the existing sheets.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:G2")
$dataRange = $ws.Range("A3:G11")
$allRange = $ws.Range("A1:G11")

$allRange.Borders.LineStyle = 1
$allRange.Font.Name = "Arial"
$allRange.Font.Size = 13

# Row 3 (bold header with fill) + Row4 header styling already mirrors the
# Ham1/PTB2 sheets: done via direct cell formatting below.
$ws.Range("A3:G3").HorizontalAlignment = -4108
$ws.Range("A3:G3").VerticalAlignment = -4108
$ws.Range("A3:G3").Interior.ColorIndex = 57
$ws.Range("A3:G3").Font.Bold = $true
$ws.Range("A3:G3").Font.Color = -16777208
$ws.Range("G3").WrapText = $true

$ws.Range("A4:G4").VerticalAlignment = -4108
$ws.Range("A4:G4").WrapText = $true

$ws.Range("A5:A11").HorizontalAlignment = -4108
$ws.Range("A5:A11").VerticalAlignment = -4105
$ws.Range("A5:A11").WrapText = $true

$ws.Range("B5:G11").VerticalAlignment = -4105
$ws.Range("B5:G11").WrapText = $true
$ws.Range("C5:E7").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 6. Row heights / column widths (approximate autosized layout).
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(3).RowHeight = 33.6
$ws.Rows.Item(4).RowHeight = 33.6
$ws.Rows.Item(5).RowHeight = 84
$ws.Rows.Item(6).RowHeight = 84
$ws.Rows.Item(7).RowHeight = 87.6

$ws.Columns.Item(1).ColumnWidth = 17.2
$ws.Columns.Item(2).ColumnWidth = 23.2
$ws.Columns.Item(3).ColumnWidth = 24.5
$ws.Columns.Item(4).ColumnWidth = 29.05
$ws.Columns.Item(5).ColumnWidth = 31.5
$ws.Columns.Item(6).ColumnWidth = 10.6
$ws.Columns.Item(7).ColumnWidth = 69.2

# ---------------------------------------------------------------------
# 7. Make the new sheet active / selected with its own cell selection.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E7").Select()

# Deselect the full-sheet-selection look for PTB2 (second sheet), matching
# the author's final saved state.
$ptb2 = $wb.Worksheets.Item("PTB2")
$ptb2.Cells.Select()
$ws.Activate()
$ws.Range("E7").Select()
